$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = -2.67

# Remove row 3 entirely (id_DK_Decentral_EP / -2.592592592592593)
$ws.Range("A3:B3").EntireRow.Delete()
